# AgeRangeMapping.xlsx corrections
# - Append the underlying numeric age ranges to the four "AggregateMapping"
#   labels in column D (Kids/Preteens/Teens/Adults) so the bar-graph
#   categories are self-describing.
# - Update the sheet's view/selection state to match where the author was
#   last working (scrolled down near row 39, cell C44 selected) instead of
#   the old selection sitting on D57:D62.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Fix up the AggregateMapping labels (column D) ---------------------
# Use whole-cell Replace so we only touch exact matches of the bare label
# (avoids accidentally re-touching a cell twice if it already has the
# age range appended, and avoids partial/substring hits).
$colD = $ws.Range("D2:D62")

$colD.Replace("Kids", "Kids (2-7)", 1, 1, $false, $false, $false, $false)
$colD.Replace("Preteens", "Preteens (8-12)", 1, 1, $false, $false, $false, $false)
$colD.Replace("Teens", "Teens (13-17)", 1, 1, $false, $false, $false, $false)
$colD.Replace("Adults", "Adults (18+)", 1, 1, $false, $false, $false, $false)

# --- 2. Update the view: scroll position + selection -----------------------
# Scroll so row 39 is at the top of the pane, then select C44 (previously the
# selection was parked on D57:D62).
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
$ws.Range("C44").Select()
